# Trade #91 closed at 2026-02-17 21:18:53 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook:
#  - Summary: refreshed aggregate metrics
#  - Strategy Status: refreshed MarketMaking strategy row
#  - All Trades / MarketMaking: trade #119 (row 120 / row 87) closed out
#    with an early exit, and a brand-new open trade #152 appended.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.05   # Current Capital
$summary.Range("B4").Value = 0.84      # Total P&L $
$summary.Range("B6").Value = 119       # Total Trades
$summary.Range("B7").Value = 53        # Winning Trades
$summary.Range("B9").Value = 44.54     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.05
$status.Range("D5").Value = 86
$status.Range("E5").Value = 0.73
$status.Range("F5").Value = 1.05
$status.Range("G5").Value = 45.35

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #119 (row 120) closes out via early exit.
$allTrades.Range("G120").Value = 0.09
$allTrades.Range("H120").Value = "CLOSED"
$allTrades.Range("I120").Value = 12.5
$allTrades.Range("J120").Value = 0.01
$allTrades.Range("K120").Value = 101.05
$allTrades.Range("L120").Value = "early_exit"
$allTrades.Range("M120").Value = 0.14

# New trade #152 appended as row 153 - duplicate the last row (which is
# already text-typed for the Date/Time columns) then adjust the cells
# that differ, so the Date column stays plain text instead of being
# auto-converted to a date serial.
$allTrades.Range("A152:Q152").Copy($allTrades.Range("A153:Q153"))
$allTrades.Range("A153").Value = 152
$allTrades.Range("C153").Value = "21:18:47"
$allTrades.Range("E153").Value = "UP"
$allTrades.Range("F153").Value = 0.08
$allTrades.Range("K153").Value = 101.0403221760222

# ---------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Same trade #119 (row 87 on this sheet) closes out via early exit.
$mm.Range("G87").Value = 0.09
$mm.Range("H87").Value = "CLOSED"
$mm.Range("I87").Value = 12.5
$mm.Range("J87").Value = 0.01
$mm.Range("K87").Value = 101.05
$mm.Range("P87").Value = "early_exit"
$mm.Range("Q87").Value = 0.14

# New trade #152 appended as row 120 on this sheet too.
$mm.Range("A119:Q119").Copy($mm.Range("A120:Q120"))
$mm.Range("A120").Value = 152
$mm.Range("C120").Value = "21:18:47"
$mm.Range("E120").Value = "UP"
$mm.Range("F120").Value = 0.08
$mm.Range("K120").Value = 101.0403221760222
